$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J (10). This shifts J:O -> K:P and
# keeps each existing cell's value/style attached to its (now shifted) cell.
$ws.Columns("J:J").Insert()

# The new column J needs the same header formatting as its neighbours
# (K10, the old J10, carries the bold/filled header style).
$ws.Range("K10").Copy()
$ws.Range("J10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J10").Value = "TOT SIN IMPU"

# Match the authored column width for the new column (~12.83 chars,
# auto best-fit like its neighbours).
$ws.Columns("J:J").ColumnWidth = 12

# Move the active selection, as recorded after the edit.
[void]$ws.Range("J11").Select()
